$p = $ppt.ActivePresentation

# "changed order of slides 7 and 8" - swap the positions of slide 7 and slide 8
$movedSlide = $p.Slides.Item(7)
$movedSlide.MoveTo(8)

# Resize the large chart picture on slide 9 ("Picture 5") - width shrinks,
# height / position stay the same.
$s9 = $p.Slides.Item(9)
$pic = $s9.Shapes.Item(7)
$pic.Width = 473.7839
